$wb = $excel.ActiveWorkbook
for ($i=0; $i -lt $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i+1)
    Write-Output $ws.Name
}
